$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.751.48'
$ws.Range("E2").Value = '  -2.60%  '
$ws.Range("D3").Value = '3.553.96'
$ws.Range("E3").Value = '  -3.54%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '617.47'
$ws.Range("E5").Value = '  -7.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.86'
$ws.Range("E6").Value = '  -3.85%  '
$ws.Range("D7").Value = '3.550.11'
$ws.Range("E7").Value = '  -3.51%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -2.22%  '
$ws.Range("E10").Value = '  -3.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.92'
$ws.Range("E11").Value = '  -2.71%  '
$ws.Range("E13").Value = '  -3.21%  '
$ws.Range("D14").Value = '4.155.84'
$ws.Range("E14").Value = '  -3.52%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '32.13'
$ws.Range("E15").Value = '  -1.88%  '
$ws.Range("D16").Value = '3.545.78'
$ws.Range("E16").Value = '  -4.06%  '
$ws.Range("D17").Value = '67.817.41'
$ws.Range("E17").Value = '  -2.50%  '
$ws.Range("E18").Value = '  -0.94%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.43'
$ws.Range("E19").Value = '  -0.42%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.67'
$ws.Range("E20").Value = '  -2.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '454.29'
$ws.Range("E21").Value = '  -2.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.60'
$ws.Range("E22").Value = '  -1.51%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '77.60'
$ws.Range("E24").Value = '  -2.87%  '
$ws.Range("D25").Value = '3.698.69'
$ws.Range("E25").Value = '  -3.45%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.64'
$ws.Range("E27").Value = '  -2.40%  '
$ws.Range("E28").Value = '  -7.30%  '
$ws.Range("E29").Value = '  -6.61%  '
$ws.Range("E30").Value = '  -3.67%  '
$ws.Range("E31").Value = '  -3.21%  '
$ws.Range("E32").Value = '  -0.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '25.94'
$ws.Range("E33").Value = '  -2.93%  '
$ws.Range("E34").Value = '  -4.40%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.23'
$ws.Range("E35").Value = '  -3.66%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.158'
$ws.Range("E36").Value = '  -3.57%  '
$ws.Range("D37").Value = '3.554.69'
$ws.Range("E37").Value = '  -3.30%  '
$ws.Range("E38").Value = '  -3.85%  '
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '176.90'
$ws.Range("E41").Value = '  -0.95%  '
$ws.Range("E42").Value = '  -1.53%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.61'
$ws.Range("E43").Value = '  -7.15%  '
$ws.Range("E44").Value = '  -5.65%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.892'
$ws.Range("E45").Value = '  -4.34%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '29.17'
$ws.Range("E46").Value = '  +6.62%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '46.19'
$ws.Range("E47").Value = '  -1.77%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.62'
$ws.Range("E48").Value = '  -4.29%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.71'
$ws.Range("E49").Value = '  -1.55%  '
$ws.Range("E50").Value = '  -6.02%  '
$ws.Range("E51").Value = '  -4.37%  '
